# Update "想去人数" (column F) counts that changed in the latest scrape
# (gh-pages data refresh at commit 456a3b4).
#
# Sheet 1 = 展览 (Exhibitions)
# Sheet 2 = 演出 (Performances)
# Sheet 3 = 本地生活 (Local life) -- no changes
# Sheet 4 = 全部类型 (All types, combined view) -- re-derived, rows shifted by +1
#           after row 6 because it includes the 演出 sheet's single row.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws4 = $wb.Worksheets.Item(4)

# --- Sheet 1: 展览 ---
$ws1.Cells.Item(2, 6).Value  = 636
$ws1.Cells.Item(3, 6).Value  = 2226
$ws1.Cells.Item(4, 6).Value  = 97
$ws1.Cells.Item(5, 6).Value  = 13364
$ws1.Cells.Item(8, 6).Value  = 523
$ws1.Cells.Item(9, 6).Value  = 488
$ws1.Cells.Item(11, 6).Value = 1005
$ws1.Cells.Item(12, 6).Value = 13821
$ws1.Cells.Item(13, 6).Value = 14499
$ws1.Cells.Item(14, 6).Value = 42
$ws1.Cells.Item(15, 6).Value = 174
$ws1.Cells.Item(21, 6).Value = 8
$ws1.Cells.Item(23, 6).Value = 4
$ws1.Cells.Item(24, 6).Value = 1107
$ws1.Cells.Item(27, 6).Value = 5541
$ws1.Cells.Item(30, 6).Value = 5351
$ws1.Cells.Item(32, 6).Value = 19
$ws1.Cells.Item(33, 6).Value = 127

# --- Sheet 2: 演出 ---
$ws2.Cells.Item(2, 6).Value = 1

# --- Sheet 4: 全部类型 ---
$ws4.Cells.Item(2, 6).Value  = 636
$ws4.Cells.Item(3, 6).Value  = 2226
$ws4.Cells.Item(4, 6).Value  = 97
$ws4.Cells.Item(5, 6).Value  = 13364
$ws4.Cells.Item(7, 6).Value  = 1
$ws4.Cells.Item(9, 6).Value  = 523
$ws4.Cells.Item(10, 6).Value = 488
$ws4.Cells.Item(12, 6).Value = 1005
$ws4.Cells.Item(13, 6).Value = 13821
$ws4.Cells.Item(14, 6).Value = 14499
$ws4.Cells.Item(15, 6).Value = 42
$ws4.Cells.Item(16, 6).Value = 174
$ws4.Cells.Item(22, 6).Value = 8
$ws4.Cells.Item(24, 6).Value = 4
$ws4.Cells.Item(25, 6).Value = 1107
$ws4.Cells.Item(28, 6).Value = 5541
$ws4.Cells.Item(31, 6).Value = 5351
$ws4.Cells.Item(33, 6).Value = 19
$ws4.Cells.Item(34, 6).Value = 127
